# Updates to ICER calculation
# Refresh the underlying cost inputs (e.g. currency/price update) that feed
# the incremental cost-effectiveness ratio (ICER) calculations. Dependent
# formulas (+/-20% ranges in columns D/E/F) recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- 1.1_TumourAgnosticCosts ---
$ws1 = $wb.Worksheets.Item("1.1_TumourAgnosticCosts")

$ws1.Range("C3").Value = 127.23357437878371
$ws1.Range("C4").Value = 95.69281677470633
$ws1.Range("C5").Value = 121.81946136426386
$ws1.Range("C6").Value = 1363.4406977269464
$ws1.Range("C9").Value = 1363.4406977269464

# c_pack_treat becomes formula-driven (still resolves to 2002)
$ws1.Range("C8").Formula = "=2002"
$ws1.Range("E8").Formula = "=C8"

# --- 1.4_AdminCost ---
$ws2 = $wb.Worksheets.Item("1.4_AdminCost")

$ws2.Range("B2").Value = 126.86521516588937
$ws2.Range("B3").Value = 106.72057071073046
$ws2.Range("B5").Value = 130.06917290304321
$ws2.Range("B6").Value = 167.28496463077394
$ws2.Range("B7").Value = 127.27194513012687
$ws2.Range("B9").Value = 2.8010648480506681
